# New patch update on Build v4.22.00.159
#
# Refresh the "Employee" roster sheet (column A) with the latest pulled
# records: two brand-new QCO entries are inserted at the top of the feed
# (Russel Bergstrom, and a newer Odell Carter record), which pushes the
# existing rows down by two positions. The sheet keeps a fixed 10-row
# window, so the two oldest entries simply fall off the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee")

$nl = [char]10

# Row 1 ("ROLE GROUP : null") is untouched.

# New record pulled to the top of the feed.
$ws.Cells.Item(2, 1).Value = "32426530 - Russel Bergstrom" + $nl + "ROLE : QCO 1706020141934"

# Everything that used to start at row 2 shifts down one row...
$ws.Cells.Item(3, 1).Value = "89212114 - Mohammed Turner" + $nl + "ROLE : RTGO100 2023-11-08T12:47:56.947450800"

# ...except row 3 (Odell Carter) is itself refreshed with a newer pull.
$ws.Cells.Item(4, 1).Value = "32695715 - Odell Carter" + $nl + "ROLE : QCO 1706021777804"

# The remaining original rows (4-7) continue to shift down by two.
$ws.Cells.Item(5, 1).Value = "92457737 - Lilliana Williamson" + $nl + "ROLE : RTGO100 1701844270281"
$ws.Cells.Item(6, 1).Value = "90317880 - Lewis Mosciski" + $nl + "ROLE : RTGO100 1701844270281"
$ws.Cells.Item(7, 1).Value = "90833312 - Angelo Mueller" + $nl + "ROLE : RTGO100 1701844270281"
$ws.Cells.Item(8, 1).Value = "92970163 - Glenna Lynch" + $nl + "ROLE : RTGO100 1701853905917"

# The old row 7 / row 8 entries (Odell Carter 1705835784686 and the
# "RTGO Operator 2024-01-09" group header) are pushed out of the fixed
# 10-row window entirely; rows 9-10 stay blank as they already were.
